# padel.xlsx update — records the new match "Miguel Ángel-Juanjo vs Luis-Raquel"
# (6-0,6-0 / 0-6,0-6) played on 2025-11-13, and propagates it through the
# derived "clasificacion_auto" / "clasificacion" standings plus the
# "historial_partidos" match log (whose FECHA column is stamped with the
# latest update date for every row).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) resultados: append the new result as row 22
# ---------------------------------------------------------------------
$wsResultados = $wb.Worksheets.Item("resultados")
$wsResultados.Range("A22").Value = "Mediocre bajo"
$wsResultados.Range("B22").Value = "1ª vuelta"
$wsResultados.Range("C22").Value = "Miguel Ángel-Juanjo"
$wsResultados.Range("D22").Value = "Luis-Raquel"
$wsResultados.Range("E22").Value = "6-0,6-0"
$wsResultados.Range("F22").Value = "0-6,0-6"
$wsResultados.Range("B22").Select()

# ---------------------------------------------------------------------
# 2) clasificacion_auto: refresh the standings for the two teams that
#    played (Miguel Ángel-Juanjo -> row 9, Luis-Raquel -> row 13)
# ---------------------------------------------------------------------
$wsAuto = $wb.Worksheets.Item("clasificacion_auto")

# Miguel Ángel-Juanjo: PUNTOS, PJ, PG, SG, JG
$wsAuto.Range("D9").Value = 6
$wsAuto.Range("E9").Value = 3
$wsAuto.Range("F9").Value = 2
$wsAuto.Range("I9").Value = 4
$wsAuto.Range("K9").Value = 31

# Luis-Raquel: PJ, PP, SP, JP
$wsAuto.Range("E13").Value = 3
$wsAuto.Range("H13").Value = 3
$wsAuto.Range("J13").Value = 6
$wsAuto.Range("L13").Value = 36

# ---------------------------------------------------------------------
# 3) clasificacion: mirror the same standings update
#    (Miguel Ángel-Juanjo -> row 15, Luis-Raquel -> row 19)
# ---------------------------------------------------------------------
$wsClasificacion = $wb.Worksheets.Item("clasificacion")

$wsClasificacion.Range("D15").Value = 6
$wsClasificacion.Range("E15").Value = 3
$wsClasificacion.Range("F15").Value = 2
$wsClasificacion.Range("I15").Value = 4

$wsClasificacion.Range("E19").Value = 3
$wsClasificacion.Range("H19").Value = 3
$wsClasificacion.Range("J19").Value = 6

# ---------------------------------------------------------------------
# 4) historial_partidos: the FECHA column is re-stamped for every
#    existing row, and two new rows log the new match (one per team)
# ---------------------------------------------------------------------
$wsHistorial = $wb.Worksheets.Item("historial_partidos")

$wsHistorial.Range("A2:A41").Value = 45974

# Copy row 41's formatting (date number format) down for the two new rows
$wsHistorial.Range("A41").Copy($wsHistorial.Range("A42"))
$wsHistorial.Range("A41").Copy($wsHistorial.Range("A43"))

$wsHistorial.Range("A42").Value = 45974
$wsHistorial.Range("B42").Value = "mediocre bajo"
$wsHistorial.Range("C42").Value = "1ª vuelta"
$wsHistorial.Range("D42").Value = "Miguel Ángel-Juanjo"
$wsHistorial.Range("E42").Value = "Gana"
$wsHistorial.Range("F42").Value = 2
$wsHistorial.Range("G42").Value = 0
$wsHistorial.Range("H42").Value = 3
$wsHistorial.Range("I42").Value = 3
$wsHistorial.Range("J42").Value = 6
$wsHistorial.Range("K42").Value = 2
$wsHistorial.Range("L42").Value = 0
$wsHistorial.Range("M42").Value = 1

$wsHistorial.Range("A43").Value = 45974
$wsHistorial.Range("B43").Value = "mediocre bajo"
$wsHistorial.Range("C43").Value = "1ª vuelta"
$wsHistorial.Range("D43").Value = "Luis-Raquel"
$wsHistorial.Range("E43").Value = "Pierde"
$wsHistorial.Range("F43").Value = 0
$wsHistorial.Range("G43").Value = 2
$wsHistorial.Range("H43").Value = 0
$wsHistorial.Range("I43").Value = 3
$wsHistorial.Range("J43").Value = 0
$wsHistorial.Range("K43").Value = 0
$wsHistorial.Range("L43").Value = 0
$wsHistorial.Range("M43").Value = 3
